$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '59.115.82'
$ws.Range('E2').Value = '  +0.65%  '
$ws.Range('D3').Value = '2.519.43'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '534.77'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.84'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.52%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.564'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.17%  '
$ws.Range('D9').Value = '2.522.88'
$ws.Range('E9').Value = '  +0.33%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0991'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.66%  '
$ws.Range('E11').Value = '  +1.73%  '
$ws.Range('E12').Value = '  -1.98%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.353'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.95%  '
$ws.Range('D14').Value = '2.970.09'
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '23.28'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.98%  '
$ws.Range('D16').Value = '59.075.72'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('D18').Value = '2.538.21'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.00'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.28%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '320.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('E22').Value = '  +0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '62.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.87%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.421'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.68%  '
$ws.Range('E26').Value = '  +1.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.47%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.82'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.80'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').Value = '0.0₃0770'
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('E31').Value = '  +0.66%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.04'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.30%  '
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('E34').Value = '  -8.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.43'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.65%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.42'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.18%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.27'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.58'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.94%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.84'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.28%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.64'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.21%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '288.89'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.64%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.34'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.806'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.90%  '
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.85'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.71%  '
$ws.Range('B46').Value = 'Mantle'
$ws.Range('C46').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.597'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '124.65'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0929'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.63'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.63%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0509'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0223'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.31%  '
